# Update the cryptos price/volume table with refreshed values from the
# Wed Feb 28 03:51:26 UTC 2024 GitHub Actions data pull.
# Price-like cells are forced to Text format before assignment so Excel's
# COM layer does not reinterpret strings such as "396.10" or "0.0000108"
# as numbers (which would silently drop precision / trailing zeros).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.941.72'
$ws.Range("E2").Value = '  +1.26%  '
$ws.Range("D3").Value = '3.246.52'
$ws.Range("E3").Value = '  +0.21%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '396.10'
$ws.Range("E5").Value = '  -1.65%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '107.75'
$ws.Range("E6").Value = '  -3.29%  '
$ws.Range("E7").Value = '  +3.90%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E9").Value = '  -1.63%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.12'
$ws.Range("E10").Value = '  -0.87%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0948'
$ws.Range("E11").Value = '  +5.46%  '
$ws.Range("E12").Value = '  +2.03%  '
$ws.Range("D13").Value = '3.755.47'
$ws.Range("E13").Value = '  -0.11%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.28'
$ws.Range("E14").Value = '  +1.86%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '18.85'
$ws.Range("E15").Value = '  -1.66%  '
$ws.Range("D16").Value = '3.239.53'
$ws.Range("E16").Value = '  -0.44%  '
$ws.Range("E17").Value = '  -3.81%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.90'
$ws.Range("E18").Value = '  +1.91%  '
$ws.Range("D19").Value = '56.771.70'
$ws.Range("E19").Value = '  +0.97%  '
$ws.Range("E20").Value = '  -2.78%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0000108'
$ws.Range("E21").Value = '  +5.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.84'
$ws.Range("E22").Value = '  -2.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '291.08'
$ws.Range("E23").Value = '  -4.73%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.03'
$ws.Range("E24").Value = '  -2.16%  '
$ws.Range("E25").Value = '  -3.01%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.07'
$ws.Range("E26").Value = '  -2.24%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '27.98'
$ws.Range("E27").Value = '  -1.55%  '
$ws.Range("E28").Value = '  +0.12%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.29'
$ws.Range("E29").Value = '  -2.14%  '
$ws.Range("E30").Value = '  -2.37%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.109'
$ws.Range("E32").Value = '  -2.78%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.16'
$ws.Range("E33").Value = '  -2.24%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '40.63'
$ws.Range("E34").Value = '  +10.70%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0482'
$ws.Range("E35").Value = '  -2.31%  '
$ws.Range("E36").Value = '  +0.88%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '51.15'
$ws.Range("E37").Value = '  -0.66%  '
$ws.Range("E38").Value = '  +0.01%  '
$ws.Range("E39").Value = '  -3.08%  '
$ws.Range("E40").Value = '  -4.30%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '138.05'
$ws.Range("E41").Value = '  +4.72%  '
$ws.Range("E42").Value = '  +1.57%  '
$ws.Range("E43").Value = '  -0.16%  '
$ws.Range("B44").Value = 'NEARProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.90'
$ws.Range("E44").Value = '  -3.28%  '
$ws.Range("B45").Value = 'ARBITRUM'
$ws.Range("C45").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.86'
$ws.Range("E45").Value = '  -3.73%  '
$ws.Range("E46").Value = '  -3.50%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '22.17'
$ws.Range("E47").Value = '  -2.35%  '
$ws.Range("B48").Value = 'WEMIXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.21'
$ws.Range("E48").Value = '  +5.43%  '
$ws.Range("D49").Value = '2.145.04'
$ws.Range("E49").Value = '  -0.80%  '
$ws.Range("E50").Value = '  -5.30%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.96'
$ws.Range("E51").Value = '  -6.71%  '
